$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force each Price (column D) cell we are about to rewrite to stay as plain
# text, matching the original inline-string cell type. Without this, Excel
# auto-converts clean-looking decimals (e.g. "92.50", "2.90", "76.00") into
# numbers and silently drops significant trailing zeros. NumberFormat is set
# per-cell (a single multi-cell Range assignment is not reliably honored for
# every cell by this engine).
$priceCellRefs = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D39", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $priceCellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.764.98'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '2.475.43'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '321.08'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '92.50'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '32.99'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('D13').Value = '2.856.22'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '6.89'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').Value = '15.53'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').Value = '2.485.46'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '0.791'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '41.705.36'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').Value = '0.0₃0942'
$ws.Range('E20').Value = '  -1.33%  '
$ws.Range('D21').Value = '71.56'
$ws.Range('E21').Value = '  +0.60%  '
$ws.Range('D22').Value = '11.25'
$ws.Range('E22').Value = '  -3.38%  '
$ws.Range('D23').Value = '240.19'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('D24').Value = '2.76'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('D27').Value = '24.87'
$ws.Range('E27').Value = '  +0.57%  '
$ws.Range('E28').Value = '  -1.46%  '
$ws.Range('D29').Value = '9.73'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').Value = '36.42'
$ws.Range('E30').Value = '  +1.75%  '
$ws.Range('D31').Value = '155.03'
$ws.Range('E31').Value = '  -0.69%  '
$ws.Range('D32').Value = '5.44'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').Value = '0.0763'
$ws.Range('E34').Value = '  -0.12%  '
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').Value = '17.10'
$ws.Range('E36').Value = '  -3.05%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').Value = '0.116'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Value = '2.90'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('E42').Value = '  -3.69%  '
$ws.Range('D43').Value = '2.004.56'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('D44').Value = '0.0283'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').Value = '18.68'
$ws.Range('E45').Value = '  -1.15%  '
$ws.Range('D46').Value = '2.96'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = '9.43'
$ws.Range('E47').Value = '  +4.45%  '
$ws.Range('D48').Value = '2.732.57'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').Value = '97.57'
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('D50').Value = '76.00'
$ws.Range('E50').Value = '  +4.05%  '
$ws.Range('D51').Value = '67.14'
$ws.Range('E51').Value = '  +0.08%  '
